$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: UNRATE_pct / Unemployment Rate -- auto-updated data + news
$ws.Range("E8").Value = 4.4

# F8 holds a "Mon YYYY" text label; force text formatting first so Excel
# doesn't auto-convert it into a date serial, then restore the plain
# (unstyled) look of the sibling cells in the column.
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Dec 2025"
$ws.Range("F8").Style = $ws.Range("F7").Style

$ws.Range("G8").Value = 4.587499999999999
$ws.Range("H8").Value = 0.3000000000000007
$ws.Range("I8").Value = 0.07317073170731725
